$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values per row (bus index 0..23, rows 2..25),
# for the scenario "case with 380 kV" where the slack bus (column B) is set to 1.02 p.u.
$newValues = @{
    2 = @(1.02, 1.030094156859033, 1.033669348053059, 1.029847730271327, 1.035820363481174, 1.035237302641636, 1.036471032017796, 1.032660452963381, 1.015751317678294)
    3 = @(1.02, 1.030983058356788, 1.034322540826321, 1.030600611673587, 1.03604439206038, 1.035767822487928, 1.036933669006539, 1.033221714495673, 1.01592708143386)
    4 = @(1.02, 1.031558417772075, 1.034745258136857, 1.031088328203587, 1.036188038287384, 1.036110655197864, 1.037232400095141, 1.033584777805449, 1.016040642581566)
    5 = @(1.02, 1.031800340905186, 1.034922980837736, 1.031293494834103, 1.036248111582146, 1.036254673442431, 1.037357835686089, 1.033737382006193, 1.016088342650856)
    6 = @(1.02, 1.031840963293059, 1.034952821940798, 1.031327950828154, 1.036258179626405, 1.036278848337185, 1.037378887989497, 1.033763003269393, 1.016096349284369)
    7 = @(1.02, 1.031561650197186, 1.034747632829699, 1.031091069138896, 1.036188842230584, 1.036112580004919, 1.037234076766192, 1.033586817020236, 1.016041280113974)
    8 = @(1.02, 1.030394527127689, 1.033890084705999, 1.03010205486912, 1.035896346893445, 1.035416686630972, 1.036627511307609, 1.032850156084643, 1.01581075277257)
    9 = @(1.02, 1.028339345207363, 1.032379485807432, 1.028363579809469, 1.035370890684403, 1.034187042864507, 1.035553919578148, 1.031551271472484, 1.015403251525917)
    10 = @(1.02, 1.026970262437105, 1.031372852802211, 1.027207574095588, 1.035013876379566, 1.033365065846301, 1.034835065501006, 1.030684883938532, 1.015130745054222)
    11 = @(1.02, 1.026377695517656, 1.030937091037144, 1.026707735996931, 1.034857702284895, 1.033008629006748, 1.034523065677704, 1.030309632843343, 1.015012552612284)
    12 = @(1.02, 1.026157629220022, 1.030775248781546, 1.02652218355063, 1.034799454871648, 1.032876156125066, 1.034407066240288, 1.030170234057846, 1.014968621710609)
    13 = @(1.02, 1.026204832384885, 1.030809963608041, 1.026561980192881, 1.034811959871821, 1.032904575446754, 1.034431953420436, 1.030200136177676, 1.014978046342808)
    14 = @(1.02, 1.026359503957696, 1.030923712708243, 1.026692395927721, 1.034852892374639, 1.03299768030881, 1.034513479341654, 1.030298110371602, 1.015008921858944)
    15 = @(1.02, 1.026454807478747, 1.030993799854123, 1.026772763944734, 1.03487808080237, 1.033055035199172, 1.034563695765329, 1.030358473708712, 1.015027941456555)
    16 = @(1.02, 1.027009594620534, 1.031401775467344, 1.027240762002326, 1.035024207811605, 1.033388710617681, 1.034855756584755, 1.030709786117219, 1.015138585014307)
    17 = @(1.02, 1.027357666997669, 1.031657720162104, 1.027534518649402, 1.035115445661573, 1.033597879161868, 1.035038763584386, 1.030930129255298, 1.015207936766112)
    18 = @(1.02, 1.027560716249512, 1.031807019558234, 1.027705931300769, 1.035168510226052, 1.03371983383853, 1.035145437844149, 1.03105864199592, 1.015248369596436)
    19 = @(1.02, 1.027629954922131, 1.031857928645916, 1.027764390305511, 1.035186577911632, 1.033761408753127, 1.035181798998411, 1.031102459842152, 1.015262152933861)
    20 = @(1.02, 1.027320319580935, 1.031630258550779, 1.027502994155188, 1.035105672517209, 1.033575442499624, 1.035019135962933, 1.030906489515292, 1.015200497928858)
    21 = @(1.02, 1.026313955936715, 1.030890215895065, 1.026653988704583, 1.034840845320907, 1.032970265347178, 1.034489474985113, 1.030269259785812, 1.014999830585543)
    22 = @(1.02, 1.025681443924036, 1.030425032752765, 1.026120820592272, 1.034672964749903, 1.032589325127033, 1.034155827077728, 1.029868528763757, 1.014873495661129)
    23 = @(1.02, 1.02601672844658, 1.030671624114554, 1.026403402374089, 1.034762091360574, 1.032791310224662, 1.034332759419795, 1.03008097100682, 1.014940483935353)
    24 = @(1.02, 1.027337195195271, 1.031642667236468, 1.027517238502995, 1.035110089053052, 1.033585580815481, 1.035028005058536, 1.030917171326279, 1.015203859277565)
    25 = @(1.02, 1.028870481199821, 1.032769942078622, 1.028812499051303, 1.035507919462469, 1.034505330343648, 1.035832024520715, 1.033737382006193, 1.015508749949986)
}

$colOrder = @("B","C","D","E","I","J","K","L","N")

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $colOrder.Length; $i++) {
        $ws.Range($colOrder[$i] + $row).Value = $vals[$i]
    }
}
